$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# ---- Row 20 -> id 3494865 ----
Set-TextValue "A20" "3494865"
$ws.Range("B20").Value = "Varta Electronics CR2430 1er Bli"
$ws.Range("C20").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-cr2430-1er-bli/p/3494865"
$ws.Range("D20").Value = "1ST"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = "Varta"
Set-TextValue "H20" "4.95"
$ws.Range("I20").Value = "4.95/1ST"
$ws.Range("J20").Value = "Preis pro 1 Stück"
Set-TextValue "K20" "4.95"
$ws.Range("L20").Value = "1ST"
$ws.Range("M20").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N20").Value = "Varta Electronics CR2430 1er Bli 4.95 Schweizer Franken"

# ---- Row 21 -> id 6753557 ----
Set-TextValue "A21" "6753557"
$ws.Range("B21").Value = "Duracell Batterien PLUS AAA/LR03 4 Stück"
$ws.Range("C21").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/duracell-batterien-plus-aaalr03-4-stueck/p/6753557"
$ws.Range("D21").Value = "4ST"
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = "Duracell"
Set-TextValue "H21" "9.95"
$ws.Range("I21").Value = "2.49/1ST"
$ws.Range("J21").Value = "Preis pro 1 Stück"
Set-TextValue "K21" "2.49"
$ws.Range("L21").Value = "1ST"
$ws.Range("M21").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']"
$ws.Range("N21").Value = "Duracell Batterien PLUS AAA/LR03 4 Stück 9.95 Schweizer Franken"

# ---- Row 22 -> id 6761133 ----
Set-TextValue "A22" "6761133"
$ws.Range("B22").Value = "Duracell Batterien PLUS C/LR14 2 Stück"
$ws.Range("C22").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterien-plus-clr14-2-stueck/p/6761133"
$ws.Range("D22").Value = "2ST"
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = "Duracell"
Set-TextValue "H22" "9.95"
$ws.Range("I22").Value = "4.98/1ST"
$ws.Range("J22").Value = "Preis pro 1 Stück"
Set-TextValue "K22" "4.98"
$ws.Range("L22").Value = "1ST"
$ws.Range("M22").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N22").Value = "Duracell Batterien PLUS C/LR14 2 Stück 9.95 Schweizer Franken"

# ---- Row 23 -> id 6801740 ----
Set-TextValue "A23" "6801740"
$ws.Range("B23").Value = "Duracell Batterien PLUS AA/LR6 12 Stück"
$ws.Range("C23").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/duracell-batterien-plus-aalr6-12-stueck/p/6801740"
$ws.Range("D23").Value = "12ST"
$ws.Range("E23").Value = ""
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = "Duracell"
Set-TextValue "H23" "29.85"
$ws.Range("I23").Value = "2.49/1ST"
$ws.Range("J23").Value = "Preis pro 1 Stück"
Set-TextValue "K23" "2.49"
$ws.Range("L23").Value = "1ST"
$ws.Range("M23").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Range("N23").Value = "Duracell Batterien PLUS AA/LR6 12 Stück 29.85 Schweizer Franken"

# ---- Row 24 -> id 3494233 ----
Set-TextValue "A24" "3494233"
$ws.Range("B24").Value = "Varta Electronics CR2032 1er Bli"
$ws.Range("C24").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-cr2032-1er-bli/p/3494233"
$ws.Range("D24").Value = "1ST"
$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 4.5
$ws.Range("G24").Value = "Varta"
Set-TextValue "H24" "4.95"
$ws.Range("I24").Value = "4.95/1ST"
$ws.Range("J24").Value = "Preis pro 1 Stück"
Set-TextValue "K24" "4.95"
$ws.Range("L24").Value = "1ST"
$ws.Range("M24").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N24").Value = "Varta Electronics CR2032 1er Bli 4.95 Schweizer Franken"

# ---- Row 26 -> id 6761135 ----
Set-TextValue "A26" "6761135"
$ws.Range("B26").Value = "Duracell Batterie PLUS 9V/6LR61 1 Stück"
$ws.Range("C26").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterie-plus-9v6lr61-1-stueck/p/6761135"
$ws.Range("D26").Value = "1ST"
$ws.Range("E26").Value = ""
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = "Duracell"
Set-TextValue "H26" "9.95"
$ws.Range("I26").Value = "9.95/1ST"
$ws.Range("J26").Value = "Preis pro 1 Stück"
Set-TextValue "K26" "9.95"
$ws.Range("L26").Value = "1ST"
$ws.Range("M26").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N26").Value = "Duracell Batterie PLUS 9V/6LR61 1 Stück 9.95 Schweizer Franken"

# ---- Row 27 -> id 6577801 ----
Set-TextValue "A27" "6577801"
$ws.Range("B27").Value = "Duracell Knopfzelle CR2025 2 Stück"
$ws.Range("C27").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-knopfzelle-cr2025-2-stueck/p/6577801"
$ws.Range("D27").Value = "2ST"
$ws.Range("E27").Value = ""
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = "Duracell"
Set-TextValue "H27" "9.95"
$ws.Range("I27").Value = "4.98/1ST"
$ws.Range("J27").Value = "Preis pro 1 Stück"
Set-TextValue "K27" "4.98"
$ws.Range("L27").Value = "1ST"
$ws.Range("M27").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N27").Value = "Duracell Knopfzelle CR2025 2 Stück 9.95 Schweizer Franken"

# ---- Row 29 -> id 6801782 ----
Set-TextValue "A29" "6801782"
$ws.Range("B29").Value = "Duracell Batterien PLUS AAA/LR03 12 Stück"
$ws.Range("C29").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/duracell-batterien-plus-aaalr03-12-stueck/p/6801782"
$ws.Range("D29").Value = "12ST"
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = "Duracell"
Set-TextValue "H29" "29.85"
$ws.Range("I29").Value = "2.49/1ST"
$ws.Range("J29").Value = "Preis pro 1 Stück"
Set-TextValue "K29" "2.49"
$ws.Range("L29").Value = "1ST"
$ws.Range("M29").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']"
$ws.Range("N29").Value = "Duracell Batterien PLUS AAA/LR03 12 Stück 29.85 Schweizer Franken"

# ---- Row 30 -> id 4014527 ----
Set-TextValue "A30" "4014527"
$ws.Range("B30").Value = "Varta Longlife Batterien AA/LR6 10 Stück"
$ws.Range("C30").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-batterien-aalr6-10-stueck/p/4014527"
$ws.Range("D30").Value = "10ST"
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = "Varta"
Set-TextValue "H30" "15.95"
$ws.Range("I30").Value = "1.60/1ST"
$ws.Range("J30").Value = "Preis pro 1 Stück"
Set-TextValue "K30" "1.60"
$ws.Range("L30").Value = "1ST"
$ws.Range("M30").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Range("N30").Value = "Varta Longlife Batterien AA/LR6 10 Stück 15.95 Schweizer Franken"

# ---- Row 31 -> id 6753554 ----
Set-TextValue "A31" "6753554"
$ws.Range("B31").Value = "Duracell Batterien Optimum AA/LR6 4 Stück"
$ws.Range("C31").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/duracell-batterien-optimum-aalr6-4-stueck/p/6753554"
$ws.Range("D31").Value = "4ST"
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 5
$ws.Range("G31").Value = "Duracell"
Set-TextValue "H31" "11.95"
$ws.Range("I31").Value = "2.99/1ST"
$ws.Range("J31").Value = "Preis pro 1 Stück"
Set-TextValue "K31" "2.99"
$ws.Range("L31").Value = "1ST"
$ws.Range("M31").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Range("N31").Value = "Duracell Batterien Optimum AA/LR6 4 Stück 11.95 Schweizer Franken"

# ---- Row 41 -> id 3494138 ----
Set-TextValue "A41" "3494138"
$ws.Range("B41").Value = "Varta Longlife Power Batterien AAA/LR03 8 Stück"
$ws.Range("C41").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-power-batterien-aaalr03-8-stueck/p/3494138"
$ws.Range("D41").Value = "8ST"
$ws.Range("E41").Value = 1
$ws.Range("F41").Value = 4
$ws.Range("G41").Value = "Varta"
Set-TextValue "H41" "14.95"
$ws.Range("I41").Value = "1.87/1ST"
$ws.Range("J41").Value = "Preis pro 1 Stück"
Set-TextValue "K41" "1.87"
$ws.Range("L41").Value = "1ST"
$ws.Range("M41").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']"
$ws.Range("N41").Value = "Varta Longlife Power Batterien AAA/LR03 8 Stück 14.95 Schweizer Franken"

# ---- Row 42 -> id 5763068 ----
Set-TextValue "A42" "5763068"
$ws.Range("B42").Value = "satrap Aspira Sine A700 Beutelloser Zyklonstaubsauger"
$ws.Range("C42").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-aspira-sine-a700-beutelloser-zyklonstaubsauger/p/5763068"
$ws.Range("D42").Value = ""
$ws.Range("E42").Value = 3
$ws.Range("F42").Value = 3.5
$ws.Range("G42").Value = "satrap"
Set-TextValue "H42" "119.00"
$ws.Range("I42").Value = ""
$ws.Range("J42").Value = ""
$ws.Range("K42").Value = ""
$ws.Range("L42").Value = ""
$ws.Range("M42").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Range("N42").Value = "satrap Aspira Sine A700 Beutelloser Zyklonstaubsauger 119.00 Schweizer Franken"

# ---- Row 43 -> id 4358323 ----
Set-TextValue "A43" "4358323"
$ws.Range("B43").Value = "Rayovac Hörgerätebatterien 312 6 Stück"
$ws.Range("C43").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/rayovac-hoergeraetebatterien-312-6-stueck/p/4358323"
$ws.Range("D43").Value = "6ST"
$ws.Range("E43").Value = 1
$ws.Range("F43").Value = 4
$ws.Range("G43").Value = "Rayovac"
Set-TextValue "H43" "9.95"
$ws.Range("I43").Value = "1.66/1ST"
$ws.Range("J43").Value = "Preis pro 1 Stück"
Set-TextValue "K43" "1.66"
$ws.Range("L43").Value = "1ST"
$ws.Range("M43").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N43").Value = "Rayovac Hörgerätebatterien 312 6 Stück 9.95 Schweizer Franken"

# ---- Row 44 -> id 5750424 ----
Set-TextValue "A44" "5750424"
$ws.Range("B44").Value = "satrap Vento HT2 Reisehaartrockner"
$ws.Range("C44").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/beautygeraete/satrap-vento-ht2-reisehaartrockner/p/5750424"
$ws.Range("D44").Value = ""
$ws.Range("E44").Value = 1
$ws.Range("F44").Value = 5
$ws.Range("G44").Value = "satrap"
Set-TextValue "H44" "29.95"
$ws.Range("I44").Value = ""
$ws.Range("J44").Value = ""
$ws.Range("K44").Value = ""
$ws.Range("L44").Value = ""
$ws.Range("M44").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'beautygeraete']"
$ws.Range("N44").Value = "satrap Vento HT2 Reisehaartrockner 29.95 Schweizer Franken"

# ---- Row 45 -> id 3494131 ----
Set-TextValue "A45" "3494131"
$ws.Range("B45").Value = "Varta Longlife Power AAA 4er Bli"
$ws.Range("C45").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-power-aaa-4er-bli/p/3494131"
$ws.Range("D45").Value = "4ST"
$ws.Range("E45").Value = 3
$ws.Range("F45").Value = 3.5
$ws.Range("G45").Value = "Varta"
Set-TextValue "H45" "8.95"
$ws.Range("I45").Value = "2.24/1ST"
$ws.Range("J45").Value = "Preis pro 1 Stück"
Set-TextValue "K45" "2.24"
$ws.Range("L45").Value = "1ST"
$ws.Range("M45").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aaa']"
$ws.Range("N45").Value = "Varta Longlife Power AAA 4er Bli 8.95 Schweizer Franken"

# ---- Row 48 -> id 5751576 ----
Set-TextValue "A48" "5751576"
$ws.Range("B48").Value = "satrap Toasty 1 Toaster"
$ws.Range("C48").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-toasty-1-toaster/p/5751576"
$ws.Range("D48").Value = ""
$ws.Range("E48").Value = 1
$ws.Range("F48").Value = 5
$ws.Range("G48").Value = "satrap"
Set-TextValue "H48" "29.95"
$ws.Range("I48").Value = ""
$ws.Range("J48").Value = ""
$ws.Range("K48").Value = ""
$ws.Range("L48").Value = ""
$ws.Range("M48").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Range("N48").Value = "satrap Toasty 1 Toaster 29.95 Schweizer Franken"

# ---- Row 49 -> id 6973029 ----
Set-TextValue "A49" "6973029"
$ws.Range("B49").Value = "Severin Tischgrill PG 8565"
$ws.Range("C49").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/severin-tischgrill-pg-8565/p/6973029"
$ws.Range("D49").Value = ""
$ws.Range("E49").Value = 1
$ws.Range("F49").Value = 3
$ws.Range("G49").Value = "Severin"
Set-TextValue "H49" "99.95"
$ws.Range("I49").Value = ""
$ws.Range("J49").Value = ""
$ws.Range("K49").Value = ""
$ws.Range("L49").Value = ""
$ws.Range("M49").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Range("N49").Value = "Severin Tischgrill PG 8565 99.95 Schweizer Franken"

# ---- Row 52 -> id 5872158 ----
Set-TextValue "A52" "5872158"
$ws.Range("B52").Value = "satrap Tischventilator Venti 1"
$ws.Range("C52").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-tischventilator-venti-1/p/5872158"
$ws.Range("D52").Value = ""
$ws.Range("E52").Value = ""
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = "satrap"
Set-TextValue "H52" "29.95"
$ws.Range("I52").Value = ""
$ws.Range("J52").Value = ""
$ws.Range("K52").Value = ""
$ws.Range("L52").Value = ""
$ws.Range("M52").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Range("N52").Value = "satrap Tischventilator Venti 1 29.95 Schweizer Franken"

# ---- Row 53 -> id 4096751 ----
Set-TextValue "A53" "4096751"
$ws.Range("B53").Value = "Varta Longlife Power Batterien AA/LR6 6 Stück"
$ws.Range("C53").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-power-batterien-aalr6-6-stueck/p/4096751"
$ws.Range("D53").Value = "6ST"
$ws.Range("E53").Value = ""
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = "Varta"
Set-TextValue "H53" "12.95"
$ws.Range("I53").Value = "2.16/1ST"
$ws.Range("J53").Value = "Preis pro 1 Stück"
Set-TextValue "K53" "2.16"
$ws.Range("L53").Value = "1ST"
$ws.Range("M53").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Range("N53").Value = "Varta Longlife Power Batterien AA/LR6 6 Stück 12.95 Schweizer Franken"

# ---- Row 60 -> id 3494007 ----
Set-TextValue "A60" "3494007"
$ws.Range("B60").Value = "Varta Longlife Power D 2er Bli"
$ws.Range("C60").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-power-d-2er-bli/p/3494007"
$ws.Range("D60").Value = "2ST"
$ws.Range("E60").Value = ""
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = "Varta"
Set-TextValue "H60" "8.95"
$ws.Range("I60").Value = "4.48/1ST"
$ws.Range("J60").Value = "Preis pro 1 Stück"
Set-TextValue "K60" "4.48"
$ws.Range("L60").Value = "1ST"
$ws.Range("M60").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N60").Value = "Varta Longlife Power D 2er Bli 8.95 Schweizer Franken"

# ---- Row 61 -> id 4119046 ----
Set-TextValue "A61" "4119046"
$ws.Range("B61").Value = "Varta Ultra Lithium AA 4er Bli"
$ws.Range("C61").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-ultra-lithium-aa-4er-bli/p/4119046"
$ws.Range("D61").Value = "4ST"
$ws.Range("E61").Value = 1
$ws.Range("F61").Value = 5
$ws.Range("G61").Value = "Varta"
Set-TextValue "H61" "14.95"
$ws.Range("I61").Value = "3.74/1ST"
$ws.Range("J61").Value = "Preis pro 1 Stück"
Set-TextValue "K61" "3.74"
$ws.Range("L61").Value = "1ST"
$ws.Range("M61").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Range("N61").Value = "Varta Ultra Lithium AA 4er Bli 14.95 Schweizer Franken"

# ---- Row 62 -> id 4905486 ----
Set-TextValue "A62" "4905486"
$ws.Range("B62").Value = "Alkaline Batterie 3LR12/4.5V"
$ws.Range("C62").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/alkaline-batterie-3lr1245v/p/4905486"
$ws.Range("D62").Value = "1ST"
$ws.Range("E62").Value = 1
$ws.Range("F62").Value = 1
$ws.Range("G62").Value = "Coop"
Set-TextValue "H62" "5.95"
$ws.Range("I62").Value = "5.95/1ST"
$ws.Range("J62").Value = "Preis pro 1 Stück"
Set-TextValue "K62" "5.95"
$ws.Range("L62").Value = "1ST"
$ws.Range("M62").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N62").Value = "Alkaline Batterie 3LR12/4.5V 5.95 Schweizer Franken"

# ---- Row 63 -> id 5831402 ----
Set-TextValue "A63" "5831402"
$ws.Range("B63").Value = "Satrap Mikrowelle Micro M2"
$ws.Range("C63").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-mikrowelle-micro-m2/p/5831402"
$ws.Range("D63").Value = ""
$ws.Range("E63").Value = ""
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = "satrap"
Set-TextValue "H63" "49.95"
$ws.Range("I63").Value = ""
$ws.Range("J63").Value = ""
$ws.Range("K63").Value = ""
$ws.Range("L63").Value = ""
$ws.Range("M63").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Range("N63").Value = "Satrap Mikrowelle Micro M2 50% Aktion 49.95 Schweizer Franken statt 99.90 Schweizer Franken"

# ---- Row 67 -> id 3494230 ----
Set-TextValue "A67" "3494230"
$ws.Range("B67").Value = "Varta Electronics V13GS / V357 1er Bli"
$ws.Range("C67").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-v13gs--v357-1er-bli/p/3494230"
$ws.Range("D67").Value = "1ST"
$ws.Range("E67").Value = ""
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = "Varta"
Set-TextValue "H67" "4.95"
$ws.Range("I67").Value = "4.95/1ST"
$ws.Range("J67").Value = "Preis pro 1 Stück"
Set-TextValue "K67" "4.95"
$ws.Range("L67").Value = "1ST"
$ws.Range("M67").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N67").Value = "Varta Electronics V13GS / V357 1er Bli 4.95 Schweizer Franken"

# ---- Row 68 -> id 4589934 ----
Set-TextValue "A68" "4589934"
$ws.Range("B68").Value = "Varta Longlife AA 4er Bli"
$ws.Range("C68").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-aa-4er-bli/p/4589934"
$ws.Range("D68").Value = "4ST"
$ws.Range("E68").Value = ""
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = "Varta"
Set-TextValue "H68" "7.95"
$ws.Range("I68").Value = "1.99/1ST"
$ws.Range("J68").Value = "Preis pro 1 Stück"
Set-TextValue "K68" "1.99"
$ws.Range("L68").Value = "1ST"
$ws.Range("M68").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Range("N68").Value = "Varta Longlife AA 4er Bli 7.95 Schweizer Franken"

# ---- Row 69 -> id 4096628 ----
Set-TextValue "A69" "4096628"
$ws.Range("B69").Value = "Varta Longlife Power Batterien AAA/LR03 6 Stück"
$ws.Range("C69").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-power-batterien-aaalr03-6-stueck/p/4096628"
$ws.Range("D69").Value = "6ST"
$ws.Range("E69").Value = ""
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = "Varta"
Set-TextValue "H69" "12.95"
$ws.Range("I69").Value = "2.16/1ST"
$ws.Range("J69").Value = "Preis pro 1 Stück"
Set-TextValue "K69" "2.16"
$ws.Range("L69").Value = "1ST"
$ws.Range("M69").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N69").Value = "Varta Longlife Power Batterien AAA/LR03 6 Stück 12.95 Schweizer Franken"

# ---- Update timestamp column O for all data rows (2-92) ----
$ws.Range("O2:O92").Value = "2022-08-11 20:57:00"
